$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix typo in the shared-string header used by column B: "Catelog" -> "Catalog"
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Catalog"

# ---------------------------------------------------------------------
# 2. Column widths were trimmed slightly (template header tweak).
#    Target character widths (from the authoritative XML) are not always
#    exactly reachable because Excel quantizes ColumnWidth to whole
#    pixels; the values below are the closest achievable inputs.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth  = 14.0
$ws.Columns.Item(2).ColumnWidth  = 12.666666666666666
$ws.Columns.Item(3).ColumnWidth  = 5.0
$ws.Columns.Item(4).ColumnWidth  = 6.833333333333333
$ws.Columns.Item(5).ColumnWidth  = 6.166666666666667
$ws.Columns.Item(6).ColumnWidth  = 12.166666666666666
$ws.Columns.Item(7).ColumnWidth  = 55.166666666666664
$ws.Columns.Item(8).ColumnWidth  = 17.0
$ws.Columns.Item(9).ColumnWidth  = 10.833333333333334
$ws.Columns.Item(10).ColumnWidth = 17.0
$ws.Columns.Item(11).ColumnWidth = 17.333333333333332
$ws.Columns.Item(12).ColumnWidth = 16.333333333333332
$ws.Columns.Item(13).ColumnWidth = 15.833333333333334
$ws.Columns.Item(14).ColumnWidth = 7.833333333333333

# ---------------------------------------------------------------------
# 3. Update the saved cursor/selection position to C7
# ---------------------------------------------------------------------
$ws.Range("C7").Select()
